# Daily attendance processing - reorder "Recorded By" (column G) values
# so that any "System" entry is moved to the end of the comma-separated
# list, preserving the relative order of the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Value2

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ',\s*'

    $others = @()
    $systemCount = 0
    foreach ($part in $parts) {
        if ($part.Equals('System')) {
            $systemCount++
        } else {
            $others += $part
        }
    }

    if ($systemCount -eq 0) {
        continue
    }

    $newParts = @()
    $newParts += $others
    for ($i = 0; $i -lt $systemCount; $i++) {
        $newParts += 'System'
    }

    $newText = $newParts -join ', '

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
